$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.811335325241089
$ws.Range("B1").Value = 1.991052031517029
$ws.Range("C1").Value = 2.349579334259033
$ws.Range("D1").Value = 2.941033363342285
$ws.Range("E1").Value = 3.38599157333374
